$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Cells.Item(2, 2).Value = 6.839249999999998
$ws.Cells.Item(2, 3).Value = 2.677449999999999
$ws.Cells.Item(2, 4).Value = 30.92719999999999
$ws.Cells.Item(2, 6).Value = 25.4711
$ws.Cells.Item(2, 7).Value = 25.4711
$ws.Cells.Item(2, 11).Value = 33.4516
$ws.Cells.Item(2, 12).Value = 25.471
$ws.Cells.Item(2, 13).Value = 7.9806
$ws.Cells.Item(2, 14).Value = 7.9806
$ws.Cells.Item(3, 2).Value = 16.919
$ws.Cells.Item(3, 3).Value = 4.097
$ws.Cells.Item(3, 4).Value = 51.64
$ws.Cells.Item(3, 6).Value = 28.217
$ws.Cells.Item(3, 7).Value = 27.166
$ws.Cells.Item(3, 8).Value = 1.05
$ws.Cells.Item(3, 9).Value = 1.05
$ws.Cells.Item(3, 11).Value = 28.1494
$ws.Cells.Item(3, 12).Value = 27.167
$ws.Cells.Item(3, 13).Value = 0.9823999999999999
$ws.Cells.Item(3, 14).Value = 0.9823999999999999
$ws.Cells.Item(4, 2).Value = 33.213
$ws.Cells.Item(4, 3).Value = 9.318
$ws.Cells.Item(4, 4).Value = 39.359
$ws.Cells.Item(4, 6).Value = 36.229
$ws.Cells.Item(4, 7).Value = 36.22893150684931
$ws.Cells.Item(4, 11).Value = 29.101
$ws.Cells.Item(4, 12).Value = 29.101

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Cells.Item(2, 2).Value = 6.166850000000001
$ws.Cells.Item(2, 3).Value = 2.487
$ws.Cells.Item(2, 4).Value = 30.87284999999999
$ws.Cells.Item(2, 6).Value = 25.3809
$ws.Cells.Item(2, 7).Value = 25.3809
$ws.Cells.Item(2, 11).Value = 185.9264
$ws.Cells.Item(2, 12).Value = 25.381
$ws.Cells.Item(2, 13).Value = 160.5454
$ws.Cells.Item(2, 14).Value = 9.0854
$ws.Cells.Item(2, 15).Value = 151.46
$ws.Cells.Item(3, 2).Value = 14.045
$ws.Cells.Item(3, 3).Value = 3.349
$ws.Cells.Item(3, 4).Value = 53.613
$ws.Cells.Item(3, 6).Value = 27.318
$ws.Cells.Item(3, 7).Value = 26.829
$ws.Cells.Item(3, 8).Value = 0.489
$ws.Cells.Item(3, 9).Value = 0.489
$ws.Cells.Item(3, 11).Value = 32.6484
$ws.Cells.Item(3, 12).Value = 26.829
$ws.Cells.Item(3, 13).Value = 5.8194
$ws.Cells.Item(3, 14).Value = 5.8194
$ws.Cells.Item(4, 2).Value = 19.833
$ws.Cells.Item(4, 3).Value = 4.341
$ws.Cells.Item(4, 4).Value = 50.695
$ws.Cells.Item(4, 6).Value = 31.69
$ws.Cells.Item(4, 7).Value = 31.68977260273973
$ws.Cells.Item(4, 11).Value = 30.9632
$ws.Cells.Item(4, 12).Value = 27.506
$ws.Cells.Item(4, 13).Value = 3.4572
$ws.Cells.Item(4, 14).Value = 3.4572

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Cells.Item(2, 2).Value = 7.666199999999999
$ws.Cells.Item(2, 3).Value = 3.051099999999999
$ws.Cells.Item(2, 4).Value = 28.0064
$ws.Cells.Item(2, 6).Value = 25.52699999999999
$ws.Cells.Item(2, 7).Value = 25.52699999999999
$ws.Cells.Item(2, 11).Value = 32.088
$ws.Cells.Item(2, 12).Value = 25.527
$ws.Cells.Item(2, 13).Value = 6.561
$ws.Cells.Item(2, 14).Value = 6.561
$ws.Cells.Item(3, 2).Value = 19.795
$ws.Cells.Item(3, 3).Value = 5.553
$ws.Cells.Item(3, 4).Value = 48.213
$ws.Cells.Item(3, 6).Value = 28.37
$ws.Cells.Item(3, 7).Value = 27.507
$ws.Cells.Item(3, 11).Value = 27.507
$ws.Cells.Item(3, 12).Value = 27.507
$ws.Cells.Item(4, 2).Value = 33.213
$ws.Cells.Item(4, 3).Value = 9.318
$ws.Cells.Item(4, 4).Value = 39.359
$ws.Cells.Item(4, 6).Value = 36.229
$ws.Cells.Item(4, 7).Value = 36.22893150684931
$ws.Cells.Item(4, 11).Value = 29.101
$ws.Cells.Item(4, 12).Value = 29.101

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Cells.Item(2, 2).Value = 8.2424
$ws.Cells.Item(2, 3).Value = 2.6334
$ws.Cells.Item(2, 4).Value = 29.6254
$ws.Cells.Item(2, 6).Value = 25.60984999999999
$ws.Cells.Item(2, 7).Value = 25.60984999999999
$ws.Cells.Item(2, 11).Value = 28.2338
$ws.Cells.Item(2, 12).Value = 25.61
$ws.Cells.Item(2, 13).Value = 2.6238
$ws.Cells.Item(2, 14).Value = 2.6238
$ws.Cells.Item(3, 2).Value = 19.795
$ws.Cells.Item(3, 3).Value = 5.553
$ws.Cells.Item(3, 4).Value = 48.213
$ws.Cells.Item(3, 6).Value = 28.37
$ws.Cells.Item(3, 7).Value = 27.507
$ws.Cells.Item(3, 11).Value = 27.507
$ws.Cells.Item(3, 12).Value = 27.507
$ws.Cells.Item(4, 2).Value = 33.213
$ws.Cells.Item(4, 3).Value = 9.318
$ws.Cells.Item(4, 4).Value = 39.359
$ws.Cells.Item(4, 6).Value = 36.229
$ws.Cells.Item(4, 7).Value = 36.22893150684931
$ws.Cells.Item(4, 11).Value = 29.101
$ws.Cells.Item(4, 12).Value = 29.101

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Cells.Item(2, 2).Value = 7.979699999999999
$ws.Cells.Item(2, 3).Value = 2.54125
$ws.Cells.Item(2, 4).Value = 29.97015
$ws.Cells.Item(2, 6).Value = 25.58135
$ws.Cells.Item(2, 7).Value = 25.58135
$ws.Cells.Item(2, 11).Value = 63.89019999999999
$ws.Cells.Item(2, 12).Value = 25.581
$ws.Cells.Item(2, 13).Value = 38.3092
$ws.Cells.Item(2, 14).Value = 3.4746
$ws.Cells.Item(2, 15).Value = 34.8346
$ws.Cells.Item(3, 2).Value = 19.795
$ws.Cells.Item(3, 3).Value = 5.553
$ws.Cells.Item(3, 4).Value = 46.576
$ws.Cells.Item(3, 6).Value = 28.337
$ws.Cells.Item(3, 7).Value = 27.473
$ws.Cells.Item(3, 8).Value = 0.864
$ws.Cells.Item(3, 9).Value = 0.864
$ws.Cells.Item(3, 11).Value = 27.6134
$ws.Cells.Item(3, 12).Value = 27.473
$ws.Cells.Item(3, 13).Value = 0.1404
$ws.Cells.Item(3, 14).Value = 0.1404
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(4, 2).Value = 33.213
$ws.Cells.Item(4, 3).Value = 9.318
$ws.Cells.Item(4, 4).Value = 38.85
$ws.Cells.Item(4, 6).Value = 36.218
$ws.Cells.Item(4, 7).Value = 36.21846575342466
$ws.Cells.Item(4, 11).Value = 29.13460000000001
$ws.Cells.Item(4, 12).Value = 29.091
$ws.Cells.Item(4, 13).Value = 0.0436
$ws.Cells.Item(4, 14).Value = 0.0436
